$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: (empty) -> "false" (as literal text, not a Boolean).
# A direct string assignment of "false"/"true" gets auto-coerced to a
# Boolean by Excel, so build it via a TEXT() formula and then collapse
# the formula down to its resulting literal value in place.
$ws.Range("B7").Formula = "=TEXT(""false"",""@"")"
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# Date: updated timestamp
$ws.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# Case Sensitive: (empty) -> "true" (as literal text, not a Boolean)
$ws.Range("B15").Formula = "=TEXT(""true"",""@"")"
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial(-4163)

$excel.CutCopyMode = $false
